# Update to version 1.3.8.3 - Append new observation log rows to the
# "15-08-2025" sheet (data/observaciones.xlsx), extending the used range
# from A1:E3 to A1:F6. Three new rows are logged by the app; the last one
# (row 6) has its columns shifted one to the right, matching the source
# data exactly (A6 holds the numeric id 11 instead of a time value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("15-08-2025")

# Copy the formatting (style index 2: left/top aligned, wrapped text)
# from the existing last data row (row 3) onto the new rows first, using
# ranges with matching column counts so no stray formatted-but-empty
# cells are introduced (rows 4 and 5 only use columns A-E).
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

$ws.Range("A3:E3").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$ws.Range("A3:E3").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# --- Row 4: new observation entry -----------------------------------
$ws.Range("A4").Value = "19:38:38"
$ws.Range("B4").Value = "LÍNEA 1"
$ws.Range("C4").Value = "MÁQUINA 31 T12"
$ws.Range("D4").Value = "[General] Comentario"
$ws.Range("E4").Value = "admin"

# --- Row 5: new observation entry -----------------------------------
$ws.Range("A5").Value = "19:43:18"
$ws.Range("B5").Value = "LÍNEA 3"
$ws.Range("C5").Value = "MÁQUINA 33 T16"
$ws.Range("D5").Value = "[General] hola"
$ws.Range("E5").Value = "admin"

# --- Row 6: new observation entry (columns shifted right by one) ----
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "19:56:23"
$ws.Range("C6").Value = "LÍNEA 2"
$ws.Range("D6").Value = "MÁQUINA 32 T23"
$ws.Range("E6").Value = "[General] aa"
$ws.Range("F6").Value = "admin"
